$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header swap between D1 (MAKE_TEXT -> MAKE) and E1 (MAKE -> MAKE_TEXT)
$ws.Range("D1").Value = "MAKE"
$ws.Range("E1").Value = "MAKE_TEXT"

# Data row updates
$ws.Range("B2").Value = "SYMBOL_2000"
$ws.Range("C2").Value = 2007

# Set column E width (closest value the engine's 1/6-char quantization allows to 20.77734375)
$ws.Columns.Item(5).ColumnWidth = 20

# Update selection to F2
$ws.Range("F2").Select()
